# Natmi following Dr Hou advice:
# Re-run of the Bmp4-Bmpr1b ligand-receptor edge table with a third
# sending/target cluster ("ECs") added alongside the existing "FAPs" and
# "sCs" clusters, and refreshed statistics for every cluster pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp4"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.504599666666667
$ws.Range("H2").Value = 13.513799
$ws.Range("I2").Value = 0.2870666703033801
$ws.Range("J2").Value = 0.2870666703033801
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.510190333333334
$ws.Range("N2").Value = 10.530571
$ws.Range("O2").Value = 0.8478537661184122
$ws.Range("P2").Value = 0.8478537661184122
$ws.Range("Q2").Value = 15.81200220546989
$ws.Range("R2").Value = 142.308019849229
$ws.Range("S2").Value = 0.2433905575437934
$ws.Range("T2").Value = 0.2433905575437934

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp4"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.504599666666667
$ws.Range("H3").Value = 13.513799
$ws.Range("I3").Value = 0.2870666703033801
$ws.Range("J3").Value = 0.2870666703033801
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.629899
$ws.Range("N3").Value = 1.889697
$ws.Range("O3").Value = 0.1521462338815877
$ws.Range("P3").Value = 0.1521462338815877
$ws.Range("Q3").Value = 2.837442825433667
$ws.Range("R3").Value = 25.536985428903
$ws.Range("S3").Value = 0.0436761127595867
$ws.Range("T3").Value = 0.04367611275958671

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp4"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.480073
$ws.Range("H4").Value = 22.440219
$ws.Range("I4").Value = 0.476686011772755
$ws.Range("J4").Value = 0.476686011772755
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.510190333333334
$ws.Range("N4").Value = 10.530571
$ws.Range("O4").Value = 0.8478537661184122
$ws.Range("P4").Value = 0.8478537661184122
$ws.Range("Q4").Value = 26.25647993722767
$ws.Range("R4").Value = 236.308319435049
$ws.Range("S4").Value = 0.4041600303374961
$ws.Range("T4").Value = 0.4041600303374961

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp4"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.480073
$ws.Range("H5").Value = 22.440219
$ws.Range("I5").Value = 0.476686011772755
$ws.Range("J5").Value = 0.476686011772755
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.629899
$ws.Range("N5").Value = 1.889697
$ws.Range("O5").Value = 0.1521462338815877
$ws.Range("P5").Value = 0.1521462338815877
$ws.Range("Q5").Value = 4.711690502627
$ws.Range("R5").Value = 42.405214523643
$ws.Range("S5").Value = 0.07252598143525886
$ws.Range("T5").Value = 0.07252598143525887

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bmp4"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.707151333333333
$ws.Range("H6").Value = 11.121454
$ws.Range("I6").Value = 0.2362473179238649
$ws.Range("J6").Value = 0.2362473179238649
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.510190333333334
$ws.Range("N6").Value = 10.530571
$ws.Range("O6").Value = 0.8478537661184122
$ws.Range("P6").Value = 0.8478537661184122
$ws.Range("Q6").Value = 13.01280677447045
$ws.Range("R6").Value = 117.115260970234
$ws.Range("S6").Value = 0.2003031782371227
$ws.Range("T6").Value = 0.2003031782371227

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bmp4"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.707151333333333
$ws.Range("H7").Value = 11.121454
$ws.Range("I7").Value = 0.2362473179238649
$ws.Range("J7").Value = 0.2362473179238649
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.629899
$ws.Range("N7").Value = 1.889697
$ws.Range("O7").Value = 0.1521462338815877
$ws.Range("P7").Value = 0.1521462338815877
$ws.Range("Q7").Value = 2.335130917715333
$ws.Range("R7").Value = 21.016178259438
$ws.Range("S7").Value = 0.03594413968674216
$ws.Range("T7").Value = 0.03594413968674216
